# Apply the two textual edits described by the diff.

$d = $word.ActiveDocument
$dash = [char]8211

# 1) The byline "GUSTAVO " / "OLIVEIRA -" / " " (three separate runs with
#    identical formatting) should end up as a single run whose text reads
#    "GUSTAVO OLIVEIRA - " (trailing space before the bold "RA 01242070").
#    Re-writing the whole phrase through Find/Replace merges the
#    identically-formatted runs into one, exactly as required.
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Replacement.ClearFormatting()
$phrase1 = "GUSTAVO OLIVEIRA " + $dash + " "
$find1.Execute(
    $phrase1,
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    $phrase1,
    2
) | Out-Null

# 2) Replace the document title "SOLICITAÇÃO DE MUDANÇAS" with
#    "GESTÃO DE MUDANÇA".
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(
    "SOLICITAÇÃO DE MUDANÇAS",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "GESTÃO DE MUDANÇA",
    2
) | Out-Null
